$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.538.58"
$ws.Range("E2").Value = "  +3.89%  "

$ws.Range("D3").Value = "3.486.53"
$ws.Range("E3").Value = "  +2.44%  "

$ws.Range("E4").Value = "  +0.02%  "

$ws.Range("E5").Value = "  +3.64%  "

$ws.Range("D6").Value = "'168.73"
$ws.Range("E6").Value = "  +3.68%  "

$ws.Range("E7").Value = "  -0.05%  "

$ws.Range("D8").Value = "3.483.45"
$ws.Range("E8").Value = "  +2.39%  "

$ws.Range("D9").Value = "'0.592"
$ws.Range("E9").Value = "  +8.08%  "

$ws.Range("D10").Value = "'7.32"
$ws.Range("E10").Value = "  +0.17%  "

$ws.Range("E11").Value = "  +6.50%  "

$ws.Range("E12").Value = "  +3.90%  "

$ws.Range("D13").Value = "4.087.83"
$ws.Range("E13").Value = "  +2.39%  "

$ws.Range("E14").Value = "  -0.24%  "

$ws.Range("E15").Value = "  +4.69%  "

$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "66.558.67"
$ws.Range("E16").Value = "  +3.91%  "

$ws.Range("B17").Value = "ShibaInu"
$ws.Range("C17").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D17").Value = "'0.0000177"
$ws.Range("E17").Value = "  +3.19%  "

$ws.Range("D18").Value = "3.487.27"
$ws.Range("E18").Value = "  +2.78%  "

$ws.Range("E19").Value = "  +3.15%  "

$ws.Range("D20").Value = "'13.99"
$ws.Range("E20").Value = "  +3.92%  "

$ws.Range("D21").Value = "'391.26"
$ws.Range("E21").Value = "  +5.02%  "

$ws.Range("D22").Value = "'7.89"
$ws.Range("E22").Value = "  +1.40%  "

$ws.Range("D23").Value = "'72.91"
$ws.Range("E23").Value = "  +3.96%  "

$ws.Range("E24").Value = "  +0.10%  "

$ws.Range("E25").Value = "  +4.68%  "

$ws.Range("E26").Value = "  +6.12%  "

$ws.Range("D27").Value = "'10.28"
$ws.Range("E27").Value = "  +8.16%  "

$ws.Range("E28").Value = "  +1.47%  "

$ws.Range("E29").Value = "  +0.11%  "

$ws.Range("D30").Value = "'6.31"
$ws.Range("E30").Value = "  +3.63%  "

$ws.Range("E31").Value = "  +5.25%  "

$ws.Range("E32").Value = "  +3.48%  "

$ws.Range("D33").Value = "'23.57"
$ws.Range("E33").Value = "  +3.60%  "

$ws.Range("E34").Value = "  +5.40%  "

$ws.Range("E35").Value = "  +0.00%  "

$ws.Range("E36").Value = "  +8.46%  "

$ws.Range("D37").Value = "'161.57"
$ws.Range("E37").Value = "  +1.25%  "

$ws.Range("D38").Value = "'0.890"
$ws.Range("E38").Value = "  +4.17%  "

$ws.Range("E39").Value = "  +5.86%  "

$ws.Range("D40").Value = "'6.79"
$ws.Range("E40").Value = "  +5.08%  "

$ws.Range("E41").Value = "  +2.66%  "

$ws.Range("D42").Value = "'26.45"
$ws.Range("E42").Value = "  +2.70%  "

$ws.Range("D43").Value = "'4.62"
$ws.Range("E43").Value = "  +6.51%  "

$ws.Range("D44").Value = "'26.80"
$ws.Range("E44").Value = "  +3.57%  "

$ws.Range("E45").Value = "  +0.93%  "

$ws.Range("D46").Value = "2.767.37"
$ws.Range("E46").Value = "  +1.58%  "

$ws.Range("E47").Value = "  +2.76%  "

$ws.Range("E48").Value = "  +3.24%  "

$ws.Range("D49").Value = "'345.67"
$ws.Range("E49").Value = "  +5.72%  "

$ws.Range("E50").Value = "  +4.37%  "

$ws.Range("D51").Value = "'33.88"
$ws.Range("E51").Value = "  +12.25%  "
